# Applies the 2025-08-15 17:36 ITA model update:
# - "solar" sheet: swap the lcoe_class rank (col P) between rows 31 and 32
# - "wind" sheet: swap cap_bnd (col M), ncap_cost~USD21_alt (col O) and
#   lcoe_class (col P) values between several adjacent row pairs
#   (rows 7/8, 26/27, 94/95), and rotate the lcoe_class (col P) ranks for
#   rows 14/15/16 and rows 24/25.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "solar" worksheet (sheet5.xml)
# ---------------------------------------------------------------------
$wsSolar = $wb.Worksheets.Item("solar")

$wsSolar.Range("P31").Value = 4
$wsSolar.Range("P32").Value = 2

# ---------------------------------------------------------------------
# "wind" worksheet (sheet6.xml)
# ---------------------------------------------------------------------
$wsWind = $wb.Worksheets.Item("wind")

# Rows 7 <-> 8 : cap_bnd / ncap_cost~USD21_alt / lcoe_class swap
$wsWind.Range("M7").Value = 0.0015
$wsWind.Range("O7").Value = 31.372355844942916
$wsWind.Range("P7").Value = 1

$wsWind.Range("M8").Value = 2.4990000000000001
$wsWind.Range("O8").Value = 39.630069093581724
$wsWind.Range("P8").Value = 3

# Rows 14/15/16 : lcoe_class rank rotation
$wsWind.Range("P14").Value = 5
$wsWind.Range("P15").Value = 4
$wsWind.Range("P16").Value = 3

# Rows 24/25 : lcoe_class rank swap
$wsWind.Range("P24").Value = 1
$wsWind.Range("P25").Value = 2

# Rows 26 <-> 27 : cap_bnd / ncap_cost~USD21_alt / lcoe_class swap
$wsWind.Range("M26").Value = 0.0045
$wsWind.Range("O26").Value = 107.03467623779731
$wsWind.Range("P26").Value = 4

$wsWind.Range("M27").Value = 3.2032500000000002
$wsWind.Range("O27").Value = 42.969340484621256
$wsWind.Range("P27").Value = 3

# Rows 94 <-> 95 : cap_bnd / ncap_cost~USD21_alt / lcoe_class swap
$wsWind.Range("M94").Value = 0.43575000000000003
$wsWind.Range("O94").Value = 120.30091898433189
$wsWind.Range("P94").Value = 5

$wsWind.Range("M95").Value = 0.03
$wsWind.Range("O95").Value = 94.983314628615091
$wsWind.Range("P95").Value = 2
